# edit.ps1 - apply the commit's changes to the active document via Word COM-interop
#
# Semantic changes being applied:
#   1. The empty paragraph right after "실습과제 1" that only held the
#      "_GoBack" bookmark loses that bookmark (becomes a truly empty paragraph).
#   2. Both occurrences of the literal password "gusgks12#$" (inside the
#      DriverManager.getConnection(...) calls) are changed to "123123".
#   3. The "_GoBack" bookmark re-appears, now anchored right after the
#      second "123123" occurrence (inside the quoted string, just before
#      the closing double-quote) - i.e. the last place the document was
#      edited, which is exactly where Word leaves _GoBack after an edit.

$d = $word.ActiveDocument

# --- Step 1: remove the pre-existing "_GoBack" bookmark -------------------
# (it currently lives alone in its own paragraph, right after "실습과제 1")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: replace the first password occurrence ------------------------
# DriverManager.getConnection("jdbc:inetdae7://...", "20165164", "gusgks12#$")
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("gusgks12#$", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Text = "123123"
}

# --- Step 3: replace the second password occurrence -----------------------
# con = DriverManager.getConnection(url, "20165164", "gusgks12#$");
# Search starting just after the first (already-replaced) occurrence so we
# land on the second one.
$rng2 = $d.Range($rng1.End, $d.Content.End)
$found2 = $rng2.Find.Execute("gusgks12#$", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $rng2.Text = "123123"

    # --- Step 4: drop the "_GoBack" bookmark right after the new text -----
    # (collapsed range sitting right before the closing quote)
    $bmRange = $d.Range($rng2.End, $rng2.End)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
